# Applies the edits described in the commit diff:
#  - Slide 1 (title slide):
#      * subtitle "Rectangle 3": shrink font, reposition/resize box, add a new
#        second line "(previously draft-gandhi-spring-twamp-srpm-07)"
#      * "Rectangle 4" (authors block): reposition to match the taller subtitle
#  - Slide 4 (History of the Draft): revision-NN -> version-NN (case-preserving)
#  - Slide 5 (Updates Since IETF-106): Revision-04 -> Version-04 in the title
#  - Slide 6 (STAMP DM message figures): "Figure 6:" -> "Figure:" (both captions)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 - title slide
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)

# "Rectangle 3" = the subtitle placeholder holding the draft name
$subtitle = $slide1.Shapes.Item(2)
$subtitle.Top = 151.93756103515625
$subtitle.Height = 59.49992370605469

$subtitleRange = $subtitle.TextFrame.TextRange

# Shrink the existing "draft-gandhi-spring-stamp-srpm-00" line to 20pt
$subtitleRange.Font.Size = 20
$oldLen = $subtitleRange.Length

# Add a new paragraph: "(" + "previously " + "draft-gandhi-spring-twamp-srpm-07)"
[void]$subtitleRange.InsertAfter("`r(")
$lenA = $subtitleRange.Length
$parenRun = $subtitleRange.Characters($oldLen + 2, 1)

[void]$subtitleRange.InsertAfter("previously ")
$lenB = $subtitleRange.Length
$prevRun = $subtitleRange.Characters($lenA + 1, $lenB - $lenA)
$prevRun.Font.Size = 20
$prevRun.Font.Italic = $false

[void]$subtitleRange.InsertAfter("draft-gandhi-spring-twamp-srpm-07)")
$lenC = $subtitleRange.Length
$draftRun = $subtitleRange.Characters($lenB + 1, $lenC - $lenB)
$draftRun.Font.Size = 20
$draftRun.Font.Italic = $true

$parenRun.Font.Size = 20
$parenRun.Font.Italic = $true

# "Rectangle 4" = authors block below the subtitle, shifts with the taller subtitle
$authors = $slide1.Shapes.Item(3)
$authors.Left = 132.0
$authors.Top = 230.0624542236328

# ---------------------------------------------------------------------------
# Slide 4 - History of the Draft: revision-NN -> version-NN
# ---------------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$history = $slide4.Shapes.Item(2).TextFrame.TextRange

$history.Paragraphs(6, 1).Runs(1, 1).Text = "Presented version-00 at IETF 104 Prague in SPRING WG"
$history.Paragraphs(8, 1).Runs(1, 1).Text = "Presented version-01 at IETF 105 Montreal in IPPM WG"
$history.Paragraphs(11, 1).Runs(1, 1).Text = "Version-02 updates included a section on stand-alone LM messages"
$history.Paragraphs(14, 1).Runs(1, 1).Text = "Presented version-04 at IETF 106 Singapore in SPRING WG"

# ---------------------------------------------------------------------------
# Slide 5 - Updates Since IETF-106 (Revision-04) -> (Version-04)
# ---------------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$title5 = $slide5.Shapes.Item(1).TextFrame.TextRange
$title5.Paragraphs(1, 1).Runs(1, 1).Text = "Updates Since IETF-106 (Version-04)"

# ---------------------------------------------------------------------------
# Slide 6 - STAMP DM message figures: "Figure 6:" -> "Figure:"
# ---------------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)

$senderFigureShape = $slide6.Shapes.Item(4)
$senderRange = $senderFigureShape.TextFrame.TextRange
$senderRange.Paragraphs(14, 1).Runs(1, 1).Text = "            Figure: Sender Control Code in STAMP DM Message"

$reflectorFigureShape = $slide6.Shapes.Item(6)
$reflectorRange = $reflectorFigureShape.TextFrame.TextRange
$reflectorRange.Paragraphs(12, 1).Runs(1, 1).Text = "            Figure: Reflector Control Code in STAMP DM Message"
